# Apply the "posId" row insertion to the "_settings" sheet (sheet2.xml).
#
# Logical change (per diff):
#   - A new row is inserted at row 9 of the paramTable (rows 9-15 shift down
#     to rows 10-16).
#   - The new row 9 holds a new parameter "posId":
#       A9 = "posId", B9 = "posId", C9 = (blank), D9 = "provided",
#       E9/F9/G9 = (blank)
#   - The paramTable (table1.xml) and the sheet dimension grow from
#     A8:G15 / A1:G15 to A8:G16 / A1:G16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_settings")

# --- 1. Insert a new row at row 9; everything below (old rows 9-15) shifts
#        down to rows 10-16, carrying its values/formulas/styles along. ---
$ws.Rows(9).Insert()

# --- 2. Grow the paramTable (ListObject) + its AutoFilter to cover the
#        extra row: A8:G15 -> A8:G16. ---
$lo = $ws.ListObjects.Item("paramTable")
$lo.Resize($ws.Range("A8:G16"))

# --- 3. Populate the new row 9 with the "posId" parameter. ---
$ws.Range("A9").Value = "posId"
$ws.Range("B9").Value = "posId"
$ws.Range("D9").Value = "provided"

# --- 4. Match formatting to the rest of the table by copying the cell
#        style from the equivalent column of a neighbouring data row
#        (values are overwritten again afterwards so only the format
#        sticks). ---
$ws.Range("A15").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B15").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("D10").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("G10").Copy()
$ws.Range("G9").PasteSpecial(-4122)

# Re-apply the text values (PasteSpecial formats-only should not have
# touched them, but make sure explicitly).
$ws.Range("A9").Value = "posId"
$ws.Range("B9").Value = "posId"
$ws.Range("D9").Value = "provided"

$excel.CutCopyMode = 0
